# ---------------------------------------------------------------------------
# PlayerPerformance_4819.xlsx maintenance edit:
#   1. ODI Batting sheet: drop the four empty placeholder cells in column B
#      (rows 2, 4, 5, 9) that used to hold blank inline strings.
#   2. Add a brand-new "ODI Batting Extra" worksheet (after "ODI Bowling")
#      with MATCH_CODE / BATTING_POSITION / NUM_4 / NUM_6 /
#      PERCENT_RUNS_OF_TOTAL / MAN_OF_MATCH columns.
# ---------------------------------------------------------------------------

function Set-TextValue($range, [string]$value) {
    # Force text storage so numeric-looking strings ("4209", "1", "0",
    # "3.24%", ...) are not silently reinterpreted as numbers/percentages.
    $range.NumberFormat = "@"
    $range.Value = $value
}

function Set-BlankText($range) {
    # Touch the cell so it is materialised (present-but-blank) rather than
    # fully absent from the sheet, matching the source export's empty
    # inline-string placeholder cells.
    $range.NumberFormat = "@"
    $range.Value = ""
}

$wb = $excel.ActiveWorkbook

# --- 1) ODI Batting: remove the 4 empty placeholder cells in column B ------
$ws2 = $wb.Worksheets.Item("ODI Batting")
$ws2.Range("B2").ClearContents()
$ws2.Range("B4").ClearContents()
$ws2.Range("B5").ClearContents()
$ws2.Range("B9").ClearContents()

# --- 2) Add the new "ODI Batting Extra" worksheet after "ODI Bowling" ------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add($null, $lastSheet)
$ws4.Name = "ODI Batting Extra"

# Header row (bold, bordered, centered/top-aligned - matches the other
# sheets' header style in this workbook)
$hdr = $ws4.Range("A1:F1")
$hdr.Font.Bold = $true
$hdr.Borders.LineStyle = 1
$hdr.HorizontalAlignment = -4108  # xlCenter
$hdr.VerticalAlignment = -4160    # xlTop

$ws4.Range("A1").Value = "MATCH_CODE"
$ws4.Range("B1").Value = "BATTING_POSITION"
$ws4.Range("C1").Value = "NUM_4"
$ws4.Range("D1").Value = "NUM_6"
$ws4.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$ws4.Range("F1").Value = "MAN_OF_MATCH"

# Data rows - MATCH_CODE column stores numeric-looking codes as text
Set-TextValue $ws4.Range("A2") "4209"
Set-BlankText $ws4.Range("B2")
Set-BlankText $ws4.Range("C2")
Set-BlankText $ws4.Range("D2")
Set-BlankText $ws4.Range("E2")
$ws4.Range("F2").Value = "NO"

Set-TextValue $ws4.Range("A3") "4210"
$ws4.Range("B3").Value = 11
Set-TextValue $ws4.Range("C3") "1"
Set-TextValue $ws4.Range("D3") "0"
Set-TextValue $ws4.Range("E3") "3.24%"
$ws4.Range("F3").Value = "NO"

Set-TextValue $ws4.Range("A4") "4211"
Set-BlankText $ws4.Range("B4")
Set-BlankText $ws4.Range("C4")
Set-BlankText $ws4.Range("D4")
Set-BlankText $ws4.Range("E4")
$ws4.Range("F4").Value = "NO"

Set-TextValue $ws4.Range("A5") "4212"
Set-BlankText $ws4.Range("B5")
Set-BlankText $ws4.Range("C5")
Set-BlankText $ws4.Range("D5")
Set-BlankText $ws4.Range("E5")
$ws4.Range("F5").Value = "NO"

Set-TextValue $ws4.Range("A6") "4660"
$ws4.Range("B6").Value = 10
Set-BlankText $ws4.Range("C6")
Set-BlankText $ws4.Range("D6")
Set-BlankText $ws4.Range("E6")
$ws4.Range("F6").Value = "NO"

Set-TextValue $ws4.Range("A7") "4666"
Set-BlankText $ws4.Range("B7")
Set-BlankText $ws4.Range("C7")
Set-BlankText $ws4.Range("D7")
Set-BlankText $ws4.Range("E7")
Set-BlankText $ws4.Range("F7")

Set-TextValue $ws4.Range("A8") "4698"
Set-BlankText $ws4.Range("B8")
Set-BlankText $ws4.Range("C8")
Set-BlankText $ws4.Range("D8")
Set-BlankText $ws4.Range("E8")
Set-BlankText $ws4.Range("F8")

Set-TextValue $ws4.Range("A9") "4699"
Set-BlankText $ws4.Range("B9")
Set-BlankText $ws4.Range("C9")
Set-BlankText $ws4.Range("D9")
Set-BlankText $ws4.Range("E9")
Set-BlankText $ws4.Range("F9")

Write-Host "edit complete"
